# Updates horarios data for Línea 141 (scrape refresh 13:41:21 -> 13:55:44)
# Applies per-cell value updates + appended rows on all three sheets,
# reflecting the re-sorted/re-scraped schedule snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet: LP1912 ---
$ws = $wb.Worksheets.Item("LP1912")

$ws.Cells.Item(2,1).Value = "Última actualización: 13:55:44"
$ws.Cells.Item(3,1).Value = "Total filas: 194"
$ws.Cells.Item(66,1).Value = "08:11:18"
$ws.Cells.Item(66,2).Value = "09:28"
$ws.Cells.Item(66,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(66,4).Value = 77
$ws.Cells.Item(66,5).Value = "LP1912"
$ws.Cells.Item(67,1).Value = "08:28:52"
$ws.Cells.Item(67,2).Value = "09:28"
$ws.Cells.Item(67,3).Value = "10_OLMOS"
$ws.Cells.Item(67,4).Value = 60
$ws.Cells.Item(67,5).Value = "LP1912"
$ws.Cells.Item(102,1).Value = "10:56:15"
$ws.Cells.Item(102,2).Value = "11:15"
$ws.Cells.Item(102,3).Value = "14_ABASTO"
$ws.Cells.Item(102,4).Value = 19
$ws.Cells.Item(102,5).Value = "LP1912"
$ws.Cells.Item(103,1).Value = "09:22:34"
$ws.Cells.Item(103,2).Value = "11:15"
$ws.Cells.Item(103,3).Value = "15X38_ABASTO"
$ws.Cells.Item(103,4).Value = 113
$ws.Cells.Item(103,5).Value = "LP1912"
$ws.Cells.Item(139,1).Value = "10:49:38"
$ws.Cells.Item(139,2).Value = "12:36"
$ws.Cells.Item(139,3).Value = "27_EL RETIRO"
$ws.Cells.Item(139,4).Value = 107
$ws.Cells.Item(139,5).Value = "LP1912"
$ws.Cells.Item(140,1).Value = "11:53:44"
$ws.Cells.Item(140,2).Value = "12:36"
$ws.Cells.Item(140,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(140,4).Value = 43
$ws.Cells.Item(140,5).Value = "LP1912"
$ws.Cells.Item(166,1).Value = "13:55:43"
$ws.Cells.Item(166,2).Value = "13:55"
$ws.Cells.Item(166,3).Value = "15_ABASTO"
$ws.Cells.Item(166,4).Value = 0
$ws.Cells.Item(166,5).Value = "LP1912"
$ws.Cells.Item(167,1).Value = "13:55:43"
$ws.Cells.Item(167,2).Value = "13:56"
$ws.Cells.Item(167,3).Value = "81_EL PELIGRO"
$ws.Cells.Item(167,4).Value = 1
$ws.Cells.Item(167,5).Value = "LP1912"
$ws.Cells.Item(168,1).Value = "13:55:43"
$ws.Cells.Item(168,2).Value = "13:58"
$ws.Cells.Item(168,3).Value = "10_OLMOS"
$ws.Cells.Item(168,4).Value = 3
$ws.Cells.Item(168,5).Value = "LP1912"
$ws.Cells.Item(169,1).Value = "13:14:31"
$ws.Cells.Item(169,2).Value = "14:02"
$ws.Cells.Item(169,3).Value = "16_SANTA ANA"
$ws.Cells.Item(169,4).Value = 48
$ws.Cells.Item(169,5).Value = "LP1912"
$ws.Cells.Item(170,1).Value = "12:33:02"
$ws.Cells.Item(170,2).Value = "14:02"
$ws.Cells.Item(170,3).Value = "10_OLMOS"
$ws.Cells.Item(170,4).Value = 89
$ws.Cells.Item(170,5).Value = "LP1912"
$ws.Cells.Item(171,1).Value = "12:46:07"
$ws.Cells.Item(171,2).Value = "14:02"
$ws.Cells.Item(171,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(171,4).Value = 76
$ws.Cells.Item(171,5).Value = "LP1912"
$ws.Cells.Item(172,1).Value = "13:14:31"
$ws.Cells.Item(172,2).Value = "14:05"
$ws.Cells.Item(172,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(172,4).Value = 51
$ws.Cells.Item(172,5).Value = "LP1912"
$ws.Cells.Item(173,1).Value = "13:41:21"
$ws.Cells.Item(173,2).Value = "14:06"
$ws.Cells.Item(173,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(173,4).Value = 25
$ws.Cells.Item(173,5).Value = "LP1912"
$ws.Cells.Item(174,1).Value = "12:46:07"
$ws.Cells.Item(174,2).Value = "14:08"
$ws.Cells.Item(174,3).Value = "16_SANTA ANA"
$ws.Cells.Item(174,4).Value = 82
$ws.Cells.Item(174,5).Value = "LP1912"
$ws.Cells.Item(175,1).Value = "12:53:26"
$ws.Cells.Item(175,2).Value = "14:09"
$ws.Cells.Item(175,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(175,4).Value = 76
$ws.Cells.Item(175,5).Value = "LP1912"
$ws.Cells.Item(176,1).Value = "13:41:21"
$ws.Cells.Item(176,2).Value = "14:14"
$ws.Cells.Item(176,3).Value = "15_ABASTO"
$ws.Cells.Item(176,4).Value = 33
$ws.Cells.Item(176,5).Value = "LP1912"
$ws.Cells.Item(177,1).Value = "12:53:26"
$ws.Cells.Item(177,2).Value = "14:16"
$ws.Cells.Item(177,3).Value = "27_EL RETIRO"
$ws.Cells.Item(177,4).Value = 83
$ws.Cells.Item(177,5).Value = "LP1912"
$ws.Cells.Item(178,1).Value = "12:33:02"
$ws.Cells.Item(178,2).Value = "14:17"
$ws.Cells.Item(178,3).Value = "27_EL RETIRO"
$ws.Cells.Item(178,4).Value = 104
$ws.Cells.Item(178,5).Value = "LP1912"
$ws.Cells.Item(179,1).Value = "12:53:26"
$ws.Cells.Item(179,2).Value = "14:17"
$ws.Cells.Item(179,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(179,4).Value = 84
$ws.Cells.Item(179,5).Value = "LP1912"
$ws.Cells.Item(180,1).Value = "12:33:02"
$ws.Cells.Item(180,2).Value = "14:18"
$ws.Cells.Item(180,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(180,4).Value = 105
$ws.Cells.Item(180,5).Value = "LP1912"
$ws.Cells.Item(181,1).Value = "12:53:26"
$ws.Cells.Item(181,2).Value = "14:27"
$ws.Cells.Item(181,3).Value = "16_SANTA ANA"
$ws.Cells.Item(181,4).Value = 94
$ws.Cells.Item(181,5).Value = "LP1912"
$ws.Cells.Item(182,1).Value = "12:33:02"
$ws.Cells.Item(182,2).Value = "14:32"
$ws.Cells.Item(182,3).Value = "14X44_ABASTO"
$ws.Cells.Item(182,4).Value = 119
$ws.Cells.Item(182,5).Value = "LP1912"
$ws.Cells.Item(183,1).Value = "13:55:43"
$ws.Cells.Item(183,2).Value = "14:33"
$ws.Cells.Item(183,3).Value = "215C_EL PATO"
$ws.Cells.Item(183,4).Value = 38
$ws.Cells.Item(183,5).Value = "LP1912"
$ws.Cells.Item(184,1).Value = "12:46:07"
$ws.Cells.Item(184,2).Value = "14:34"
$ws.Cells.Item(184,3).Value = "215C_EL PATO"
$ws.Cells.Item(184,4).Value = 108
$ws.Cells.Item(184,5).Value = "LP1912"
$ws.Cells.Item(185,1).Value = "12:46:07"
$ws.Cells.Item(185,2).Value = "14:39"
$ws.Cells.Item(185,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(185,4).Value = 113
$ws.Cells.Item(185,5).Value = "LP1912"
$ws.Cells.Item(186,1).Value = "12:53:26"
$ws.Cells.Item(186,2).Value = "14:47"
$ws.Cells.Item(186,3).Value = "215B_EL PATO"
$ws.Cells.Item(186,4).Value = 114
$ws.Cells.Item(186,5).Value = "LP1912"
$ws.Cells.Item(187,1).Value = "13:41:21"
$ws.Cells.Item(187,2).Value = "14:51"
$ws.Cells.Item(187,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(187,4).Value = 70
$ws.Cells.Item(187,5).Value = "LP1912"
$ws.Cells.Item(188,1).Value = "13:55:43"
$ws.Cells.Item(188,2).Value = "14:51"
$ws.Cells.Item(188,3).Value = "16_SANTA ANA"
$ws.Cells.Item(188,4).Value = 56
$ws.Cells.Item(188,5).Value = "LP1912"
$ws.Cells.Item(189,1).Value = "13:55:43"
$ws.Cells.Item(189,2).Value = "14:53"
$ws.Cells.Item(189,3).Value = "215_EL PELIGRO"
$ws.Cells.Item(189,4).Value = 58
$ws.Cells.Item(189,5).Value = "LP1912"
$ws.Cells.Item(190,1).Value = "13:14:31"
$ws.Cells.Item(190,2).Value = "14:54"
$ws.Cells.Item(190,3).Value = "215_EL PELIGRO"
$ws.Cells.Item(190,4).Value = 100
$ws.Cells.Item(190,5).Value = "LP1912"
$ws.Cells.Item(191,1).Value = "13:14:31"
$ws.Cells.Item(191,2).Value = "15:02"
$ws.Cells.Item(191,3).Value = "10_OLMOS"
$ws.Cells.Item(191,4).Value = 108
$ws.Cells.Item(191,5).Value = "LP1912"
$ws.Cells.Item(192,1).Value = "13:14:31"
$ws.Cells.Item(192,2).Value = "15:13"
$ws.Cells.Item(192,3).Value = "17X38_ROMERO"
$ws.Cells.Item(192,4).Value = 119
$ws.Cells.Item(192,5).Value = "LP1912"
$ws.Cells.Item(193,1).Value = "13:55:43"
$ws.Cells.Item(193,2).Value = "15:17"
$ws.Cells.Item(193,3).Value = "14_ABASTO"
$ws.Cells.Item(193,4).Value = 82
$ws.Cells.Item(193,5).Value = "LP1912"
$ws.Cells.Item(194,1).Value = "13:41:21"
$ws.Cells.Item(194,2).Value = "15:18"
$ws.Cells.Item(194,3).Value = "14_ABASTO"
$ws.Cells.Item(194,4).Value = 97
$ws.Cells.Item(194,5).Value = "LP1912"
$ws.Cells.Item(195,1).Value = "13:55:43"
$ws.Cells.Item(195,2).Value = "15:33"
$ws.Cells.Item(195,3).Value = "215C_EL PATO"
$ws.Cells.Item(195,4).Value = 98
$ws.Cells.Item(195,5).Value = "LP1912"
$ws.Cells.Item(196,1).Value = "13:41:21"
$ws.Cells.Item(196,2).Value = "15:34"
$ws.Cells.Item(196,3).Value = "215C_EL PATO"
$ws.Cells.Item(196,4).Value = 113
$ws.Cells.Item(196,5).Value = "LP1912"
$ws.Cells.Item(197,1).Value = "13:55:43"
$ws.Cells.Item(197,2).Value = "15:41"
$ws.Cells.Item(197,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(197,4).Value = 106
$ws.Cells.Item(197,5).Value = "LP1912"
$ws.Cells.Item(198,1).Value = "13:55:43"
$ws.Cells.Item(198,2).Value = "15:53"
$ws.Cells.Item(198,3).Value = "15X38_ABASTO"
$ws.Cells.Item(198,4).Value = 118
$ws.Cells.Item(198,5).Value = "LP1912"
$ws.Cells.Item(199,1).Value = "13:55:43"
$ws.Cells.Item(199,2).Value = "15:53"
$ws.Cells.Item(199,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(199,4).Value = 118
$ws.Cells.Item(199,5).Value = "LP1912"

# --- Sheet: LP1912-215 ---
$ws = $wb.Worksheets.Item("LP1912-215")

$ws.Cells.Item(2,1).Value = "Última actualización: 13:55:44"
$ws.Cells.Item(3,1).Value = "Total filas: 35"
$ws.Cells.Item(34,1).Value = "13:55:43"
$ws.Cells.Item(34,2).Value = "14:33"
$ws.Cells.Item(34,3).Value = "215C_EL PATO"
$ws.Cells.Item(34,4).Value = 38
$ws.Cells.Item(34,5).Value = "LP1912"
$ws.Cells.Item(35,1).Value = "12:46:07"
$ws.Cells.Item(35,2).Value = "14:34"
$ws.Cells.Item(35,3).Value = "215C_EL PATO"
$ws.Cells.Item(35,4).Value = 108
$ws.Cells.Item(35,5).Value = "LP1912"
$ws.Cells.Item(36,1).Value = "12:53:26"
$ws.Cells.Item(36,2).Value = "14:47"
$ws.Cells.Item(36,3).Value = "215B_EL PATO"
$ws.Cells.Item(36,4).Value = 114
$ws.Cells.Item(36,5).Value = "LP1912"
$ws.Cells.Item(37,1).Value = "13:55:43"
$ws.Cells.Item(37,2).Value = "14:53"
$ws.Cells.Item(37,3).Value = "215_EL PELIGRO"
$ws.Cells.Item(37,4).Value = 58
$ws.Cells.Item(37,5).Value = "LP1912"
$ws.Cells.Item(38,1).Value = "13:14:31"
$ws.Cells.Item(38,2).Value = "14:54"
$ws.Cells.Item(38,3).Value = "215_EL PELIGRO"
$ws.Cells.Item(38,4).Value = 100
$ws.Cells.Item(38,5).Value = "LP1912"
$ws.Cells.Item(39,1).Value = "13:55:43"
$ws.Cells.Item(39,2).Value = "15:33"
$ws.Cells.Item(39,3).Value = "215C_EL PATO"
$ws.Cells.Item(39,4).Value = 98
$ws.Cells.Item(39,5).Value = "LP1912"
$ws.Cells.Item(40,1).Value = "13:41:21"
$ws.Cells.Item(40,2).Value = "15:34"
$ws.Cells.Item(40,3).Value = "215C_EL PATO"
$ws.Cells.Item(40,4).Value = 113
$ws.Cells.Item(40,5).Value = "LP1912"

# --- Sheet: 6203-6173 ---
$ws = $wb.Worksheets.Item("6203-6173")

$ws.Cells.Item(2,1).Value = "Última actualización: 13:55:44"
$ws.Cells.Item(3,1).Value = "Total filas: 31"
$ws.Cells.Item(20,1).Value = "08:52:40"
$ws.Cells.Item(20,2).Value = "10:30"
$ws.Cells.Item(20,3).Value = "215A_LA PLATA"
$ws.Cells.Item(20,4).Value = 98
$ws.Cells.Item(20,5).Value = "L6173"
$ws.Cells.Item(21,1).Value = "08:38:24"
$ws.Cells.Item(21,2).Value = "10:30"
$ws.Cells.Item(21,3).Value = "215B_LP-P MOR-1 Y 57"
$ws.Cells.Item(21,4).Value = 112
$ws.Cells.Item(21,5).Value = "L6173"
$ws.Cells.Item(31,1).Value = "13:55:43"
$ws.Cells.Item(31,2).Value = "13:58"
$ws.Cells.Item(31,3).Value = "215C_LA PLATA"
$ws.Cells.Item(31,4).Value = 3
$ws.Cells.Item(31,5).Value = "L6203"
$ws.Cells.Item(32,1).Value = "13:14:31"
$ws.Cells.Item(32,2).Value = "14:03"
$ws.Cells.Item(32,3).Value = "215C_LA PLATA"
$ws.Cells.Item(32,4).Value = 49
$ws.Cells.Item(32,5).Value = "L6203"
$ws.Cells.Item(33,1).Value = "13:55:43"
$ws.Cells.Item(33,2).Value = "14:26"
$ws.Cells.Item(33,3).Value = "215C_LA PLATA"
$ws.Cells.Item(33,4).Value = 31
$ws.Cells.Item(33,5).Value = "L6203"
$ws.Cells.Item(34,1).Value = "12:46:07"
$ws.Cells.Item(34,2).Value = "14:27"
$ws.Cells.Item(34,3).Value = "215C_LA PLATA"
$ws.Cells.Item(34,4).Value = 101
$ws.Cells.Item(34,5).Value = "L6203"
$ws.Cells.Item(35,1).Value = "13:55:43"
$ws.Cells.Item(35,2).Value = "15:21"
$ws.Cells.Item(35,3).Value = "215A_LA PLATA"
$ws.Cells.Item(35,4).Value = 86
$ws.Cells.Item(35,5).Value = "L6173"
$ws.Cells.Item(36,1).Value = "13:41:21"
$ws.Cells.Item(36,2).Value = "15:22"
$ws.Cells.Item(36,3).Value = "215A_LA PLATA"
$ws.Cells.Item(36,4).Value = 101
$ws.Cells.Item(36,5).Value = "L6173"

